# Adds SVR (Support Vector Regression) parameter columns to the
# pred_par parameter sheet: svr_kernel_scale, svr_epsilon,
# svr_box_constraint — with default values in row 2, alongside the
# existing RNN parameters (tmax_training, tmax_pred, ... , data_type).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header labels (columns K, L, M) appended after existing
# "data_type" header in column J.
$ws.Range("K1").Value = "svr_kernel_scale"
$ws.Range("L1").Value = "svr_epsilon"
$ws.Range("M1").Value = "svr_box_constraint"

# Default parameter values for the new SVR columns, row 2 (alongside
# the other parameter defaults already present in A2:J2).
$ws.Range("K2").Value = 100
$ws.Range("L2").Value = 0.2
$ws.Range("M2").Value = 20

# Tidy up the stray formatted-but-empty cell left below the remarks
# block now that the style table has been cleaned up.
$ws.Range("A13").ClearFormats()

# Matches the author's last selection position after the edit.
$ws.Range("J11").Select()
